{"js": "// Update the answers of two-digit division worksheet table.\n// Each cell in the \"answer\" rows (0, 4, 8, 12, 16 - zero based) of the\n// single table in the document holds one equation string, e.g. \"81\u00f75=16, 1\".\n// The commit replaces 25 of these equation strings with new ones while\n// leaving every other formatting property (fonts, size, alignment, etc.)\n// untouched. We target each cell explicitly by (row, col) so that\n// duplicate text values (the same string can be both an old value in one\n// cell and a new value in another) never get double-replaced.\n\nconst edits = [\n  { row: 0, col: 0, oldText: \"81\u00f75=16, 1\", newText: \"11\u00f77=1, 4\" },\n  { row: 0, col: 1, oldText: \"18\u00f76=3, 0\", newText: \"52\u00f79=5, 7\" },\n  { row: 0, col: 2, oldText: \"46\u00f78=5, 6\", newText: \"52\u00f77=7, 3\" },\n  { row: 0, col: 3, oldText: \"80\u00f72=40, 0\", newText: \"57\u00f72=28, 1\" },\n  { row: 0, col: 4, oldText: \"70\u00f73=23, 1\", newText: \"91\u00f79=10, 1\" },\n  { row: 4, col: 0, oldText: \"14\u00f75=2, 4\", newText: \"19\u00f73=6, 1\" },\n  { row: 4, col: 1, oldText: \"65\u00f79=7, 2\", newText: \"48\u00f76=8, 0\" },\n  { row: 4, col: 2, oldText: \"21\u00f76=3, 3\", newText: \"30\u00f76=5, 0\" },\n  { row: 4, col: 3, oldText: \"86\u00f75=17, 1\", newText: \"27\u00f77=3, 6\" },\n  { row: 4, col: 4, oldText: \"14\u00f73=4, 2\", newText: \"88\u00f76=14, 4\" },\n  { row: 8, col: 0, oldText: \"99\u00f77=14, 1\", newText: \"10\u00f75=2, 0\" },\n  { row: 8, col: 1, oldText: \"66\u00f77=9, 3\", newText: \"22\u00f73=7, 1\" },\n  { row: 8, col: 2, oldText: \"55\u00f74=13, 3\", newText: \"89\u00f75=17, 4\" },\n  { row: 8, col: 3, oldText: \"10\u00f78=1, 2\", newText: \"14\u00f79=1, 5\" },\n  { row: 8, col: 4, oldText: \"10\u00f77=1, 3\", newText: \"36\u00f79=4, 0\" },\n  { row: 12, col: 0, oldText: \"11\u00f75=2, 1\", newText: \"33\u00f79=3, 6\" },\n  { row: 12, col: 1, oldText: \"10\u00f74=2, 2\", newText: \"93\u00f73=31, 0\" },\n  { row: 12, col: 2, oldText: \"15\u00f77=2, 1\", newText: \"28\u00f76=4, 4\" },\n  { row: 12, col: 3, oldText: \"83\u00f75=16, 3\", newText: \"31\u00f74=7, 3\" },\n  { row: 12, col: 4, oldText: \"11\u00f77=1, 4\", newText: \"26\u00f76=4, 2\" },\n  { row: 16, col: 0, oldText: \"73\u00f72=36, 1\", newText: \"83\u00f77=11, 6\" },\n  { row: 16, col: 1, oldText: \"76\u00f76=12, 4\", newText: \"39\u00f78=4, 7\" },\n  { row: 16, col: 2, oldText: \"45\u00f75=9, 0\", newText: \"38\u00f78=4, 6\" },\n  { row: 16, col: 3, oldText: \"70\u00f76=11, 4\", newText: \"71\u00f79=7, 8\" },\n  { row: 16, col: 4, oldText: \"11\u00f78=1, 3\", newText: \"42\u00f72=21, 0\" },\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Grab the first paragraph range of every edited cell so we can replace the\n// run text in place (this keeps the existing run/paragraph formatting -\n// rFonts, sz, jc - intact instead of falling back to default formatting).\nconst ranges = edits.map((edit) => {\n  const cell = table.getCell(edit.row, edit.col);\n  const paragraph = cell.body.paragraphs.getFirst();\n  return paragraph.getRange();\n});\nranges.forEach((range) => range.load(\"text\"));\nawait context.sync();\n\nfor (let i = 0; i < edits.length; i++) {\n  const edit = edits[i];\n  const actualText = ranges[i].text;\n  if (actualText !== edit.oldText) {\n    throw new Error(\n      `Unexpected text in cell (${edit.row}, ${edit.col}): expected \"${edit.oldText}\" but found \"${actualText}\"`\n    );\n  }\n  ranges[i].insertText(edit.newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Update the answers of two-digit division worksheet table.\n# Each cell in the \"answer\" rows (1, 5, 9, 13, 17 - one based, as COM\n# indexes rows/columns starting at 1) of the single table in the document\n# holds one equation string, e.g. \"81\u00f75=16, 1\". The commit replaces 25 of\n# these equation strings with new ones while leaving every other\n# formatting property (fonts, size, alignment, etc.) untouched. We target\n# each cell explicitly by (row, col) so that duplicate text values (the\n# same string can be both an old value in one cell and a new value in\n# another) never get double-replaced.\n\n$edits = @(\n    @{ Row = 1; Col = 1; OldText = \"81\u00f75=16, 1\"; NewText = \"11\u00f77=1, 4\" },\n    @{ Row = 1; Col = 2; OldText = \"18\u00f76=3, 0\"; NewText = \"52\u00f79=5, 7\" },\n    @{ Row = 1; Col = 3; OldText = \"46\u00f78=5, 6\"; NewText = \"52\u00f77=7, 3\" },\n    @{ Row = 1; Col = 4; OldText = \"80\u00f72=40, 0\"; NewText = \"57\u00f72=28, 1\" },\n    @{ Row = 1; Col = 5; OldText = \"70\u00f73=23, 1\"; NewText = \"91\u00f79=10, 1\" },\n    @{ Row = 5; Col = 1; OldText = \"14\u00f75=2, 4\"; NewText = \"19\u00f73=6, 1\" },\n    @{ Row = 5; Col = 2; OldText = \"65\u00f79=7, 2\"; NewText = \"48\u00f76=8, 0\" },\n    @{ Row = 5; Col = 3; OldText = \"21\u00f76=3, 3\"; NewText = \"30\u00f76=5, 0\" },\n    @{ Row = 5; Col = 4; OldText = \"86\u00f75=17, 1\"; NewText = \"27\u00f77=3, 6\" },\n    @{ Row = 5; Col = 5; OldText = \"14\u00f73=4, 2\"; NewText = \"88\u00f76=14, 4\" },\n    @{ Row = 9; Col = 1; OldText = \"99\u00f77=14, 1\"; NewText = \"10\u00f75=2, 0\" },\n    @{ Row = 9; Col = 2; OldText = \"66\u00f77=9, 3\"; NewText = \"22\u00f73=7, 1\" },\n    @{ Row = 9; Col = 3; OldText = \"55\u00f74=13, 3\"; NewText = \"89\u00f75=17, 4\" },\n    @{ Row = 9; Col = 4; OldText = \"10\u00f78=1, 2\"; NewText = \"14\u00f79=1, 5\" },\n    @{ Row = 9; Col = 5; OldText = \"10\u00f77=1, 3\"; NewText = \"36\u00f79=4, 0\" },\n    @{ Row = 13; Col = 1; OldText = \"11\u00f75=2, 1\"; NewText = \"33\u00f79=3, 6\" },\n    @{ Row = 13; Col = 2; OldText = \"10\u00f74=2, 2\"; NewText = \"93\u00f73=31, 0\" },\n    @{ Row = 13; Col = 3; OldText = \"15\u00f77=2, 1\"; NewText = \"28\u00f76=4, 4\" },\n    @{ Row = 13; Col = 4; OldText = \"83\u00f75=16, 3\"; NewText = \"31\u00f74=7, 3\" },\n    @{ Row = 13; Col = 5; OldText = \"11\u00f77=1, 4\"; NewText = \"26\u00f76=4, 2\" },\n    @{ Row = 17; Col = 1; OldText = \"73\u00f72=36, 1\"; NewText = \"83\u00f77=11, 6\" },\n    @{ Row = 17; Col = 2; OldText = \"76\u00f76=12, 4\"; NewText = \"39\u00f78=4, 7\" },\n    @{ Row = 17; Col = 3; OldText = \"45\u00f75=9, 0\"; NewText = \"38\u00f78=4, 6\" },\n    @{ Row = 17; Col = 4; OldText = \"70\u00f76=11, 4\"; NewText = \"71\u00f79=7, 8\" },\n    @{ Row = 17; Col = 5; OldText = \"11\u00f78=1, 3\"; NewText = \"42\u00f72=21, 0\" }\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\nforeach ($edit in $edits) {\n    $cell = $t.Cell($edit.Row, $edit.Col)\n    $range = $cell.Range\n    # Cell ranges include the trailing end-of-cell marker (cell mark), so\n    # trim it off before comparing / when deciding what to overwrite.\n    $currentText = $range.Text.TrimEnd([char]0x07, [char]0x0D)\n    if ($currentText -ne $edit.OldText) {\n        throw \"Unexpected text in cell ($($edit.Row), $($edit.Col)): expected '$($edit.OldText)' but found '$currentText'\"\n    }\n    $range.Text = $edit.NewText\n}\n"}
